# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-label suffixes to the new
#   version-tagged suffixes "_FV2310" / "_FV2404" (row 1 headers).
# - Turn the header+data range into a real Excel Table (adds the table
#   part + autofilter) so the sheet keeps its filter/sort affordances.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 21 header cells in row 1.
$oldSuffix = "_old"
$newSuffix = "_new"
$fv2310 = "_FV2310"
$fv2404 = "_FV2404"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value()
    if ($text.EndsWith($oldSuffix)) {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = $base + $fv2310
    } elseif ($text.EndsWith($newSuffix)) {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value = $base + $fv2404
    }
}

# 2) Convert the used range into an Excel Table (adds xl/tables/table1.xml,
#    the autoFilter and the worksheet <tableParts> relationship).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), 0, 1)
$tbl.Name = "Table1"

# 3) Freeze panes at row 1 (split under the header row).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
